$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 0.824751214671096
$ws.Cells.Item(2, 3).Value = 0.1269536563578555
$ws.Cells.Item(2, 4).Value = 0.11414641847729
$ws.Cells.Item(2, 5).Value = 0.1359348193697016
$ws.Cells.Item(2, 6).Value = 2.061461219676175
$ws.Cells.Item(2, 10).Value = 0.1915216948782703
$ws.Cells.Item(2, 11).Value = 0.4103022769977542
$ws.Cells.Item(2, 12).Value = 0.2480145936236227
$ws.Cells.Item(2, 13).Value = 0.22758329612072
$ws.Cells.Item(2, 14).Value = 2.577820989056841
$ws.Cells.Item(2, 15).Value = 5.513561611086914

$ws.Cells.Item(3, 2).Value = 0.7932138792440355
$ws.Cells.Item(3, 3).Value = 0.1258845308745791
$ws.Cells.Item(3, 4).Value = 0.1129672082035071
$ws.Cells.Item(3, 5).Value = 0.1363442329701634
$ws.Cells.Item(3, 6).Value = 2.069261271783304
$ws.Cells.Item(3, 10).Value = 0.1925720608257988
$ws.Cells.Item(3, 11).Value = 0.3805653120324166
$ws.Cells.Item(3, 12).Value = 0.2463367822808209
$ws.Cells.Item(3, 13).Value = 0.2221226912225802
$ws.Cells.Item(3, 14).Value = 2.599368906662384
$ws.Cells.Item(3, 15).Value = 5.540438367404818

$ws.Cells.Item(4, 2).Value = 0.7741695568240061
$ws.Cells.Item(4, 3).Value = 0.1252223891874209
$ws.Cells.Item(4, 4).Value = 0.1122816244530611
$ws.Cells.Item(4, 5).Value = 0.1366315677627448
$ws.Cells.Item(4, 6).Value = 2.074851106953332
$ws.Cells.Item(4, 10).Value = 0.1932567539112222
$ws.Cells.Item(4, 11).Value = 0.3624249812154403
$ws.Cells.Item(4, 12).Value = 0.245391029690829
$ws.Cells.Item(4, 13).Value = 0.218858270694863
$ws.Cells.Item(4, 14).Value = 2.613275470802918
$ws.Cells.Item(4, 15).Value = 5.559102859912599

$ws.Cells.Item(5, 2).Value = 0.7664899202037248
$ws.Cells.Item(5, 3).Value = 0.1249511436995974
$ws.Cells.Item(5, 4).Value = 0.1120119679782476
$ws.Cells.Item(5, 5).Value = 0.1367577237962756
$ws.Cells.Item(5, 6).Value = 2.077330612572062
$ws.Cells.Item(5, 10).Value = 0.1935457884093381
$ws.Cells.Item(5, 11).Value = 0.355062894433388
$ws.Cells.Item(5, 12).Value = 0.2450269495040942
$ws.Cells.Item(5, 13).Value = 0.2175503693991736
$ws.Cells.Item(5, 14).Value = 2.619112601451905
$ws.Cells.Item(5, 15).Value = 5.567253082521702

$ws.Cells.Item(6, 2).Value = 0.7652196394226962
$ws.Cells.Item(6, 3).Value = 0.1249060184149471
$ws.Cells.Item(6, 4).Value = 0.1119677806585742
$ws.Cells.Item(6, 5).Value = 0.1367792201068632
$ws.Cells.Item(6, 6).Value = 2.077754517477253
$ws.Cells.Item(6, 10).Value = 0.1935943878546382
$ws.Cells.Item(6, 11).Value = 0.3538422670270762
$ws.Cells.Item(6, 12).Value = 0.2449677847188383
$ws.Cells.Item(6, 13).Value = 0.2173345491237377
$ws.Cells.Item(6, 14).Value = 2.620092128513761
$ws.Cells.Item(6, 15).Value = 5.568639312952101

$ws.Cells.Item(7, 2).Value = 0.7740656573683111
$ws.Cells.Item(7, 3).Value = 0.125218736794487
$ws.Cells.Item(7, 4).Value = 0.1122779483256267
$ws.Cells.Item(7, 5).Value = 0.136633232412164
$ws.Cells.Item(7, 6).Value = 2.074883729765091
$ws.Cells.Item(7, 10).Value = 0.1932606113518469
$ws.Cells.Item(7, 11).Value = 0.3623255704119401
$ws.Cells.Item(7, 12).Value = 0.2453860331261737
$ws.Cells.Item(7, 13).Value = 0.218840541087026
$ws.Cells.Item(7, 14).Value = 2.61335350356496
$ws.Cells.Item(7, 15).Value = 5.559210572027496

$ws.Cells.Item(8, 2).Value = 0.8138111788691162
$ws.Cells.Item(8, 3).Value = 0.1265862121436427
$ws.Cells.Item(8, 4).Value = 0.1137318772403191
$ws.Cells.Item(8, 5).Value = 0.1360685374365929
$ws.Cells.Item(8, 6).Value = 2.063984692217481
$ws.Cells.Item(8, 10).Value = 0.1918756185718618
$ws.Cells.Item(8, 11).Value = 0.4000247124054681
$ws.Cells.Item(8, 12).Value = 0.2474186143401482
$ws.Cells.Item(8, 13).Value = 0.2256822139392369
$ws.Cells.Item(8, 14).Value = 2.585110453701308
$ws.Cells.Item(8, 15).Value = 5.522380352461454

$ws.Cells.Item(9, 2).Value = 0.8942635222008732
$ws.Cells.Item(9, 3).Value = 0.1292221245128786
$ws.Cells.Item(9, 4).Value = 0.1168859981138226
$ws.Cells.Item(9, 5).Value = 0.1352454212205547
$ws.Cells.Item(9, 6).Value = 2.048951993145032
$ws.Cells.Item(9, 10).Value = 0.189474361527374
$ws.Cells.Item(9, 11).Value = 0.4748729671521801
$ws.Cells.Item(9, 12).Value = 0.2520708131980598
$ws.Cells.Item(9, 13).Value = 0.2397947410174659
$ws.Cells.Item(9, 14).Value = 2.535085119439056
$ws.Cells.Item(9, 15).Value = 5.467285822606669

$ws.Cells.Item(10, 2).Value = 0.9548754209506285
$ws.Cells.Item(10, 3).Value = 0.1311303396923691
$ws.Cells.Item(10, 4).Value = 0.1193854585071392
$ws.Cells.Item(10, 5).Value = 0.1348126362440567
$ws.Cells.Item(10, 6).Value = 2.041758068376481
$ws.Cells.Item(10, 10).Value = 0.1879008878295405
$ws.Cells.Item(10, 11).Value = 0.5304062477878517
$ws.Cells.Item(10, 12).Value = 0.2558907208576571
$ws.Cells.Item(10, 13).Value = 0.2505814126513215
$ws.Cells.Item(10, 14).Value = 2.501590564579953
$ws.Cells.Item(10, 15).Value = 5.4372189116024

$ws.Cells.Item(11, 2).Value = 0.9827703845231781
$ws.Cells.Item(11, 3).Value = 0.1319921700570887
$ws.Cells.Item(11, 4).Value = 0.1205615504132282
$ws.Cells.Item(11, 5).Value = 0.13465281334115
$ws.Cells.Item(11, 6).Value = 2.039318502097075
$ws.Cells.Item(11, 10).Value = 0.1872262501672566
$ws.Cells.Item(11, 11).Value = 0.555783863692568
$ws.Cells.Item(11, 12).Value = 0.2577149046965985
$ws.Cells.Item(11, 13).Value = 0.2555781344225423
$ws.Cells.Item(11, 14).Value = 2.487058976693891
$ws.Cells.Item(11, 15).Value = 5.425794931977947

$ws.Cells.Item(12, 2).Value = 0.9933792229523135
$ws.Cells.Item(12, 3).Value = 0.1323176149247161
$ws.Cells.Item(12, 4).Value = 0.1210124739074843
$ws.Cells.Item(12, 5).Value = 0.1345975973789724
$ws.Cells.Item(12, 6).Value = 2.038514218647592
$ws.Cells.Item(12, 10).Value = 0.1869766811334621
$ws.Cells.Item(12, 11).Value = 0.565409841254791
$ws.Cells.Item(12, 12).Value = 0.2584180283317892
$ws.Cells.Item(12, 13).Value = 0.2574830475074634
$ws.Cells.Item(12, 14).Value = 2.481657566229508
$ws.Cells.Item(12, 15).Value = 5.421792477846736

$ws.Cells.Item(13, 2).Value = 0.9910924023363066
$ws.Cells.Item(13, 3).Value = 0.1322475653506103
$ws.Cells.Item(13, 4).Value = 0.120915112938178
$ws.Cells.Item(13, 5).Value = 0.1346092535213987
$ws.Cells.Item(13, 6).Value = 2.038682123524708
$ws.Cells.Item(13, 10).Value = 0.1870301680481692
$ws.Cells.Item(13, 11).Value = 0.5633360107524368
$ws.Cells.Item(13, 12).Value = 0.258266050504524
$ws.Cells.Item(13, 13).Value = 0.2570722251694804
$ws.Cells.Item(13, 14).Value = 2.482816347841453
$ws.Cells.Item(13, 15).Value = 5.422640095642123

$ws.Cells.Item(14, 2).Value = 0.9836422698919876
$ws.Cells.Item(14, 3).Value = 0.1320189629658586
$ws.Cells.Item(14, 4).Value = 0.1205985369893199
$ws.Cells.Item(14, 5).Value = 0.1346481644902546
$ws.Cells.Item(14, 6).Value = 2.039249939166794
$ws.Cells.Item(14, 10).Value = 0.1872055997795226
$ws.Cells.Item(14, 11).Value = 0.5565754806279983
$ws.Cells.Item(14, 12).Value = 0.2577725042205259
$ws.Cells.Item(14, 13).Value = 0.2557345978480328
$ws.Cells.Item(14, 14).Value = 2.486612567618364
$ws.Cells.Item(14, 15).Value = 5.425459165962309

$ws.Cells.Item(15, 2).Value = 0.9790847683079562
$ws.Cells.Item(15, 3).Value = 0.1318788182057347
$ws.Cells.Item(15, 4).Value = 0.1204053478900491
$ws.Cells.Item(15, 5).Value = 0.1346726888428673
$ws.Cells.Item(15, 6).Value = 2.039613301070446
$ws.Cells.Item(15, 10).Value = 0.1873138249083333
$ws.Cells.Item(15, 11).Value = 0.5524365314189481
$ws.Cells.Item(15, 12).Value = 0.2574717977943379
$ws.Cells.Item(15, 13).Value = 0.2549169201065453
$ws.Cells.Item(15, 14).Value = 2.488951064576147
$ws.Cells.Item(15, 15).Value = 5.427228049841972

$ws.Cells.Item(16, 2).Value = 0.9530588430472449
$ws.Cells.Item(16, 3).Value = 0.1310738903719013
$ws.Cells.Item(16, 4).Value = 0.1193093796412441
$ws.Cells.Item(16, 5).Value = 0.1348238247090841
$ws.Cells.Item(16, 6).Value = 2.041934239770384
$ws.Cells.Item(16, 10).Value = 0.1879458036328963
$ws.Cells.Item(16, 11).Value = 0.5287500332151467
$ws.Cells.Item(16, 12).Value = 0.2557732386468743
$ws.Cells.Item(16, 13).Value = 0.2502566593886755
$ws.Cells.Item(16, 14).Value = 2.502554427055566
$ws.Cells.Item(16, 15).Value = 5.438010800517134

$ws.Cells.Item(17, 2).Value = 0.9371748148161316
$ws.Cells.Item(17, 3).Value = 0.1305784869397328
$ws.Cells.Item(17, 4).Value = 0.1186470068350118
$ws.Cells.Item(17, 5).Value = 0.1349260156976158
$ws.Cells.Item(17, 6).Value = 2.043571226628117
$ws.Cells.Item(17, 10).Value = 0.1883440293788539
$ws.Cells.Item(17, 11).Value = 0.5142482622906925
$ws.Cells.Item(17, 12).Value = 0.2547533138551543
$ws.Cells.Item(17, 13).Value = 0.2474206351560539
$ws.Cells.Item(17, 14).Value = 2.511080260692042
$ws.Cells.Item(17, 15).Value = 5.445202510721629

$ws.Cells.Item(18, 2).Value = 0.9280691223756037
$ws.Cells.Item(18, 3).Value = 0.1302929589113191
$ws.Cells.Item(18, 4).Value = 0.1182697087642168
$ws.Cells.Item(18, 5).Value = 0.1349882822222597
$ws.Cells.Item(18, 6).Value = 2.044591197990997
$ws.Cells.Item(18, 10).Value = 0.1885769517827143
$ws.Cells.Item(18, 11).Value = 0.5059181076535992
$ws.Cells.Item(18, 12).Value = 0.2541748263212469
$ws.Cells.Item(18, 13).Value = 0.2457978915095254
$ws.Cells.Item(18, 14).Value = 2.516050477947186
$ws.Cells.Item(18, 15).Value = 5.449551142767547

$ws.Cells.Item(19, 2).Value = 0.9249913313711318
$ws.Cells.Item(19, 3).Value = 0.1301961840352703
$ws.Cells.Item(19, 4).Value = 0.1181425960814266
$ws.Cells.Item(19, 5).Value = 0.1350099645244072
$ws.Cells.Item(19, 6).Value = 2.04495001922853
$ws.Cells.Item(19, 10).Value = 0.1886564810336182
$ws.Cells.Item(19, 11).Value = 0.5030995458201346
$ws.Cells.Item(19, 12).Value = 0.2539803623138965
$ws.Cells.Item(19, 13).Value = 0.2452499164540782
$ws.Cells.Item(19, 14).Value = 2.517744710401626
$ws.Cells.Item(19, 15).Value = 5.451059967408185

$ws.Cells.Item(20, 2).Value = 0.9388625577312553
$ws.Cells.Item(20, 3).Value = 0.1306312841527273
$ws.Cells.Item(20, 4).Value = 0.1187171369602567
$ws.Cells.Item(20, 5).Value = 0.1349147763617733
$ws.Cells.Item(20, 6).Value = 2.043388852273722
$ws.Cells.Item(20, 10).Value = 0.1883012367840529
$ws.Cells.Item(20, 11).Value = 0.5157908784404128
$ws.Cells.Item(20, 12).Value = 0.2548610441545094
$ws.Cells.Item(20, 13).Value = 0.2477216600407317
$ws.Cells.Item(20, 14).Value = 2.510165800939925
$ws.Cells.Item(20, 15).Value = 5.444414987783375

$ws.Cells.Item(21, 2).Value = 0.9858293220303551
$ws.Cells.Item(21, 3).Value = 0.1320861338968982
$ws.Cells.Item(21, 4).Value = 0.1206913724865259
$ws.Cells.Item(21, 5).Value = 0.1346365915861529
$ws.Cells.Item(21, 6).Value = 2.039079916041203
$ws.Cells.Item(21, 10).Value = 0.1871539111770328
$ws.Cells.Item(21, 11).Value = 0.5585607816491347
$ws.Cells.Item(21, 12).Value = 0.2579171364276931
$ws.Cells.Item(21, 13).Value = 0.2561271462420009
$ws.Cells.Item(21, 14).Value = 2.48549477445707
$ws.Cells.Item(21, 15).Value = 5.424622359409597

$ws.Cells.Item(22, 2).Value = 1.016790395048503
$ws.Cells.Item(22, 3).Value = 0.1330316413579666
$ws.Cells.Item(22, 4).Value = 0.1220140491169133
$ws.Cells.Item(22, 5).Value = 0.1344856958021516
$ws.Cells.Item(22, 6).Value = 2.036960355956722
$ws.Cells.Item(22, 10).Value = 0.1864384598886417
$ws.Cells.Item(22, 11).Value = 0.5866065384729495
$ws.Cells.Item(22, 12).Value = 0.2599863816482468
$ws.Cells.Item(22, 13).Value = 0.2616949394524895
$ws.Cells.Item(22, 14).Value = 2.469961779551836
$ws.Cells.Item(22, 15).Value = 5.413572447943636

$ws.Cells.Item(23, 2).Value = 1.0002418273952
$ws.Cells.Item(23, 3).Value = 0.1325274984622808
$ws.Cells.Item(23, 4).Value = 0.1213051656006598
$ws.Cells.Item(23, 5).Value = 0.134563410564505
$ws.Cells.Item(23, 6).Value = 2.038027952215515
$ws.Cells.Item(23, 10).Value = 0.1868171676516521
$ws.Cells.Item(23, 11).Value = 0.5716296477589822
$ws.Cells.Item(23, 12).Value = 0.2588754369873243
$ws.Cells.Item(23, 13).Value = 0.2587165530826141
$ws.Cells.Item(23, 14).Value = 2.478197968089344
$ws.Cells.Item(23, 15).Value = 5.419297617743354

$ws.Cells.Item(24, 2).Value = 0.9380994473480655
$ws.Cells.Item(24, 3).Value = 0.1306074167562414
$ws.Cells.Item(24, 4).Value = 0.1186854201966554
$ws.Cells.Item(24, 5).Value = 0.1349198467134443
$ws.Cells.Item(24, 6).Value = 2.043471058101062
$ws.Cells.Item(24, 10).Value = 0.1883205709268054
$ws.Cells.Item(24, 11).Value = 0.5150934395017259
$ws.Cells.Item(24, 12).Value = 0.2548123147319359
$ws.Cells.Item(24, 13).Value = 0.2475855426182818
$ws.Cells.Item(24, 14).Value = 2.510579014498393
$ws.Cells.Item(24, 15).Value = 5.444770360198845

$ws.Cells.Item(25, 2).Value = 0.8722330020125639
$ws.Cells.Item(25, 3).Value = 0.1285139878659223
$ws.Cells.Item(25, 4).Value = 0.1160005557353543
$ws.Cells.Item(25, 5).Value = 0.1354378026512606
$ws.Cells.Item(25, 6).Value = 2.052341478158581
$ws.Cells.Item(25, 10).Value = 0.1900903902015791
$ws.Cells.Item(25, 11).Value = 0.4545279432996381
$ws.Cells.Item(25, 12).Value = 0.2507413779843191
$ws.Cells.Item(25, 13).Value = 0.2549169201065453
$ws.Cells.Item(25, 14).Value = 2.488951064576147
$ws.Cells.Item(25, 15).Value = 5.427228049841972

